# Update cryptocurrency price/volume data (scheduled GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.037.84'
$ws.Range('E2').Value = '  -2.10%  '

$ws.Range('D3').Value = '1.667.97'
$ws.Range('E3').Value = '  -1.44%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '''216.99'
$ws.Range('E5').Value = '  -1.27%  '

$ws.Range('D6').Value = '''0.5103'
$ws.Range('E6').Value = '  -0.21%  '

$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('D8').Value = '''0.2661'
$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').Value = '''0.06401'
$ws.Range('E9').Value = '  +1.20%  '

$ws.Range('D10').Value = '''21.79'
$ws.Range('E10').Value = '  -1.25%  '

$ws.Range('D11').Value = '''0.07458'
$ws.Range('E11').Value = '  +1.26%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.672.57'
$ws.Range('E12').Value = '  -1.18%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.515'
$ws.Range('E13').Value = '  -0.12%  '

$ws.Range('E14').Value = '  +0.84%  '

$ws.Range('D15').Value = '''0.000008558'
$ws.Range('E15').Value = '  +0.53%  '

$ws.Range('D16').Value = '''64.31'
$ws.Range('E16').Value = '  -1.65%  '

$ws.Range('D17').Value = '26.099.10'
$ws.Range('E17').Value = '  -1.93%  '

$ws.Range('D18').Value = '''4.941'
$ws.Range('E18').Value = '  -0.83%  '

$ws.Range('D19').Value = '''1.005'
$ws.Range('E19').Value = '  -0.07%  '

$ws.Range('E20').Value = '  -1.55%  '

$ws.Range('D21').Value = '''191.97'
$ws.Range('E21').Value = '  +2.74%  '

$ws.Range('D22').Value = '''6.196'
$ws.Range('E22').Value = '  -1.01%  '

$ws.Range('D23').Value = '''1.006'
$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('D24').Value = '''144.81'
$ws.Range('E24').Value = '  -0.05%  '

$ws.Range('D25').Value = '''7.610'
$ws.Range('E25').Value = '  +1.59%  '

$ws.Range('D26').Value = '''0.1199'
$ws.Range('E26').Value = '  +2.13%  '

$ws.Range('E27').Value = '  -0.93%  '

$ws.Range('D28').Value = '''0.06523'
$ws.Range('E28').Value = '  +13.51%  '

$ws.Range('D29').Value = '''1.338'
$ws.Range('E29').Value = '  -0.28%  '

$ws.Range('D30').Value = '''1.315'
$ws.Range('E30').Value = '  -1.77%  '

$ws.Range('D31').Value = '''3.540'
$ws.Range('E31').Value = '  +0.46%  '

$ws.Range('D32').Value = '''3.517'
$ws.Range('E32').Value = '  +0.24%  '

$ws.Range('D33').Value = '''1.652'
$ws.Range('E33').Value = '  +0.55%  '

$ws.Range('E34').Value = '  -0.03%  '

$ws.Range('D35').Value = '''0.6120'
$ws.Range('E35').Value = '  +2.07%  '

$ws.Range('E36').Value = '  +0.33%  '

$ws.Range('D37').Value = '''2.683'
$ws.Range('E37').Value = '  +0.05%  '

$ws.Range('D38').Value = '''6.262'
$ws.Range('E38').Value = '  +7.41%  '

$ws.Range('D39').Value = '''0.01601'
$ws.Range('E39').Value = '  -1.21%  '

$ws.Range('D40').Value = '1.092.45'
$ws.Range('E40').Value = '  +0.13%  '

$ws.Range('D41').Value = '''0.8693'
$ws.Range('E41').Value = '  +0.89%  '

$ws.Range('E42').Value = '  +0.32%  '

$ws.Range('D43').Value = '''101.16'
$ws.Range('E43').Value = '  +1.64%  '

$ws.Range('D44').Value = '1.816.70'
$ws.Range('E44').Value = '  -1.81%  '

$ws.Range('E45').Value = '  -3.06%  '

$ws.Range('D46').Value = '''56.46'
$ws.Range('E46').Value = '  -0.01%  '

$ws.Range('D47').Value = '''1.005'
$ws.Range('E47').Value = '  +0.04%  '

$ws.Range('D48').Value = '''8.079'
$ws.Range('E48').Value = '  -0.27%  '

$ws.Range('D49').Value = '''0.05231'
$ws.Range('E49').Value = '  -0.05%  '

$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.4287'
$ws.Range('E50').Value = '  -0.83%  '

$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '''6.077'
$ws.Range('E51').Value = '  +4.90%  '
